# LOB1232.xlsx edit: restructure rows 12-24 of the single worksheet.
#
# Strategy: many target cells reuse text that already exists elsewhere on
# the sheet (rows shift up by one, a couple of values get duplicated, one
# old row is removed and one brand new value "Semestral" is introduced).
# To avoid Excel's automatic text->number/date reinterpretation when
# round-tripping plain strings (e.g. "01/01/2018"), we stage the required
# source cells (with exact value + number format) into a scratch row far
# below the used range using Range.Copy (which preserves the literal
# stored value/type), then copy them from the scratch row into their final
# destinations, and finally remove the scratch row again.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Step 1: stage copies of every source cell we will need, into row 200.
# ---------------------------------------------------------------------
$ws.Range("B8").Copy($ws.Range("A200"))
$ws.Range("C8").Copy($ws.Range("B200"))
$ws.Range("B13").Copy($ws.Range("C200"))
$ws.Range("C13").Copy($ws.Range("D200"))
$ws.Range("A14").Copy($ws.Range("E200"))
$ws.Range("A15").Copy($ws.Range("F200"))
$ws.Range("B15").Copy($ws.Range("G200"))
$ws.Range("C15").Copy($ws.Range("H200"))
$ws.Range("A16").Copy($ws.Range("I200"))
$ws.Range("A17").Copy($ws.Range("J200"))
$ws.Range("B17").Copy($ws.Range("K200"))
$ws.Range("C17").Copy($ws.Range("L200"))
$ws.Range("A18").Copy($ws.Range("M200"))
$ws.Range("A19").Copy($ws.Range("N200"))
$ws.Range("A20").Copy($ws.Range("O200"))
$ws.Range("A21").Copy($ws.Range("P200"))
$ws.Range("A22").Copy($ws.Range("Q200"))
$ws.Range("A23").Copy($ws.Range("R200"))
$ws.Range("B24").Copy($ws.Range("S200"))
$ws.Range("C24").Copy($ws.Range("T200"))

# ---------------------------------------------------------------------
# Step 2: clear the cells that must end up empty. Use Clear() (not
# Delete(), which shifts subsequent rows up/left) and not ClearContents()
# (which would leave a bare, styled <c> element behind).
# ---------------------------------------------------------------------
$ws.Range("B17:C17").Clear()
$ws.Range("B22:C22").Clear()

# ---------------------------------------------------------------------
# Step 3: write the new / moved values from the scratch row into their
# final destinations (order-independent now, since everything we need
# was already captured in step 1).
# ---------------------------------------------------------------------
$ws.Range("B13").Copy($ws.Range("B10"))
$ws.Range("C13").Copy($ws.Range("C10"))

$ws.Range("A14").Copy($ws.Range("A13"))
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

$ws.Range("A15").Copy($ws.Range("A14"))
$ws.Range("G200").Copy($ws.Range("B14"))
$ws.Range("H200").Copy($ws.Range("C14"))

$ws.Range("A16").Copy($ws.Range("A15"))
$ws.Range("A200").Copy($ws.Range("B15"))
$ws.Range("B200").Copy($ws.Range("C15"))

$ws.Range("A17").Copy($ws.Range("A16"))
$ws.Range("K200").Copy($ws.Range("B16"))
$ws.Range("L200").Copy($ws.Range("C16"))

$ws.Range("A18").Copy($ws.Range("A17"))

$ws.Range("A19").Copy($ws.Range("A18"))
$ws.Range("C200").Copy($ws.Range("B18"))
$ws.Range("D200").Copy($ws.Range("C18"))

$ws.Range("A20").Copy($ws.Range("A19"))
$ws.Range("A21").Copy($ws.Range("A20"))
$ws.Range("A22").Copy($ws.Range("A21"))
$ws.Range("A23").Copy($ws.Range("A22"))
$ws.Range("A23").Clear()

$ws.Range("S200").Copy($ws.Range("B23"))
$ws.Range("T200").Copy($ws.Range("C23"))

# ---------------------------------------------------------------------
# Step 4: drop the now-superseded old row 24 and remove the scratch row.
# ---------------------------------------------------------------------
$ws.Rows.Item(200).Delete()
$ws.Rows.Item(24).Delete()

# ---------------------------------------------------------------------
# Step 5: row heights for the shifted rows.
# ---------------------------------------------------------------------
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(14).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(16).RowHeight = 120
$ws.Rows.Item(17).AutoFit()
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(20).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
$ws.Rows.Item(22).AutoFit()
$ws.Rows.Item(23).RowHeight = 30
